# Daily attendance processing - 2025-12-18 13:46:25
# Swap the order of names in the "Recorded By" (column G) cells that list
# "dnasr281@gmail.com" first together with exactly one other recorder,
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".
# Cells that only contain "dnasr281@gmail.com" (no other recorder) or that
# list more than two recorders are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    $parts = $text.Split(",")

    if ($parts.Count -eq 2) {
        $first = $parts[0].Trim()
        $second = $parts[1].Trim()

        if ($first -eq "dnasr281@gmail.com") {
            $cell.Value = "$second, $first"
        }
    }
}
